$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version bump: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date update
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now set
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail"; becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row; remove it entirely,
# shifting all following rows up by one.
$ws.Rows.Item(11).Delete()

# After the shift, "Case Sensitive" (now row 14) gets a value of "true"
$ws.Range("B14").Value = "true"
